$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the bordered / centered / bold style once on a scratch cell, then
# copy the finished format onto the real target cells. This avoids leaving
# behind unreferenced intermediate cellXf entries that piecemeal property
# assignment on multiple distinct cells would otherwise accumulate.
$tmp = $ws.Range("D10")
$tmp.Font.Bold = $true
$tmp.Borders.LineStyle = 1
$tmp.HorizontalAlignment = -4108  # xlCenter
$tmp.VerticalAlignment = -4160    # xlTop

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0

$tmp.Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$tmp.Clear()

# B2 holds the plain (unstyled) label text.
$ws.Range("B2").Value = "disconnected_elements"
